$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last datalist entry (recording62) by clearing row 10's contents.
# This also drops the now-unused shared string from xl/sharedStrings.xml.
$ws.Range("A10:B10").ClearContents()

# Set column A to a fixed custom width of 69 (matches width="69.0" in the
# saved worksheet XML; Excel's ColumnWidth->XML width conversion adds ~5/6
# of a character width, so 68.15 characters serializes to an XML width of 69).
$ws.Columns.Item(1).ColumnWidth = 68.15
